# Apply dynamic backtesting / model-selection re-run updates to the
# auditoria_validacion.xlsx workbook values (backtest_ranges,
# signals_distribution, dm_tests sheets).

$wb = $excel.ActiveWorkbook

# --- Sheet: backtest_ranges ---
$ws1 = $wb.Worksheets.Item("backtest_ranges")

$ws1.Range("F2").Value = 0.003734280000000006
$ws1.Range("G2").Value = 0.00500295932424001
$ws1.Range("H2").Value = 0.338699475190435
$ws1.Range("L2").Value = 0.003734280000000006
$ws1.Range("M2").Value = 0.00500295932424001
$ws1.Range("N2").Value = 0.338699475190435

$ws1.Range("F3").Value = 0.004453308288650507
$ws1.Range("G3").Value = 0.006011904051877777
$ws1.Range("H3").Value = 0.4041872880502484
$ws1.Range("L3").Value = 0.004453308288650507
$ws1.Range("M3").Value = 0.006011904051877777
$ws1.Range("N3").Value = 0.4041872880502484

$ws1.Range("F4").Value = 0.01580155341018982
$ws1.Range("G4").Value = 0.02099472482408699
$ws1.Range("H4").Value = 1.436601101621157
$ws1.Range("L4").Value = 0.01580155341018982
$ws1.Range("M4").Value = 0.02099472482408699
$ws1.Range("N4").Value = 1.436601101621157

# --- Sheet: signals_distribution ---
$ws2 = $wb.Worksheets.Item("signals_distribution")

$ws2.Range("D5").Value = 51.4
$ws2.Range("F5").Value = 48.6

$ws2.Range("D7").Value = 38.2
$ws2.Range("E7").Value = 30
$ws2.Range("F7").Value = 31.8

# --- Sheet: dm_tests ---
$ws3 = $wb.Worksheets.Item("dm_tests")

$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 1

$ws3.Range("C3").Value = [double]"-8.210391889330084E-16"
$ws3.Range("D3").Value = 0.9999999999999993

$ws3.Range("C4").Value = [double]"2.49418492728443E-15"
$ws3.Range("D4").Value = 0.999999999999998
